$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'264.85"
$ws.Range("E2").Value = "'1.52%"

# Row 3
$ws.Range("D3").Value = "'26.59"
$ws.Range("E3").Value = "'-1.65%"

# Row 4
$ws.Range("D4").Value = "'4.704"
$ws.Range("E4").Value = "'-0.04%"

# Row 5
$ws.Range("D5").Value = "'0.06097"
$ws.Range("E5").Value = "'-1.34%"

# Row 6
$ws.Range("D6").Value = "'6.734"
$ws.Range("E6").Value = "'0.77%"

# Row 7
$ws.Range("E7").Value = "'0.09%"

# Row 8
$ws.Range("D8").Value = "'0.9077"
$ws.Range("E8").Value = "'-0.50%"

# Row 9
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1406"
$ws.Range("E9").Value = "'0.02%"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.05041"
$ws.Range("E10").Value = "'8.05%"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07098"
$ws.Range("E11").Value = "'0.22%"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03142"
$ws.Range("E12").Value = "'1.00%"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09026"
$ws.Range("E13").Value = "'-0.20%"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001539"
$ws.Range("E14").Value = "'0.72%"

# Row 15
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006053"
$ws.Range("E15").Value = "'-1.81%"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006011"
$ws.Range("E16").Value = "'-0.80%"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.451"
$ws.Range("E17").Value = "'-0.05%"

# Row 18
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.169"
$ws.Range("E18").Value = "'0.16%"

# Row 19
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.175"
$ws.Range("E19").Value = "'-0.16%"

# Row 20
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3071"
$ws.Range("E20").Value = "'-0.23%"

# Row 21
$ws.Range("D21").Value = "'0.1281"
$ws.Range("E21").Value = "'-1.41%"

# Row 22
$ws.Range("D22").Value = "'4.113"
$ws.Range("E22").Value = "'0.70%"

# Row 23
$ws.Range("E23").Value = "'0.54%"

# Row 24
$ws.Range("D24").Value = "'0.001179"
$ws.Range("E24").Value = "'-3.00%"

# Row 25
$ws.Range("D25").Value = "'0.004060"
$ws.Range("E25").Value = "'6.78%"

# Row 27
$ws.Range("D27").Value = "'0.0001682"
$ws.Range("E27").Value = "'6.59%"

# Row 40
$ws.Range("D40").Value = "'0.03927"
$ws.Range("E40").Value = "'1.23%"

# Row 41
$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'0.24%"

# Row 42
$ws.Range("D42").Value = "'0.004198"
$ws.Range("E42").Value = "'2.65%"

# Row 43
$ws.Range("D43").Value = "'0.002106"
$ws.Range("E43").Value = "'-3.55%"

# Row 44
$ws.Range("D44").Value = "'0.01160"
$ws.Range("E44").Value = "'-28.97%"

# Row 45
$ws.Range("D45").Value = "'0.00005104"
$ws.Range("E45").Value = "'-1.07%"

# Row 46
$ws.Range("E46").Value = "'0.05%"

# Row 48
$ws.Range("D48").Value = "'0.2579"
$ws.Range("E48").Value = "'53.37%"

# Row 49
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.05%"

# Row 50
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.05%"
